$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "Overview"; New = "1 Overview" },
    @{ Old = "When to Use Articles vs Vignettes"; New = "2 When to Use Articles vs Vignettes" },
    @{ Old = "Advanced Quarto Features"; New = "3 Advanced Quarto Features" },
    @{ Old = "Cross-References"; New = "3.1 Cross-References" },
    @{ Old = "Summary Tables"; New = "3.2 Summary Tables" },
    @{ Old = "Code Folding"; New = "3.3 Code Folding" },
    @{ Old = "Tabsets"; New = "3.4 Tabsets" },
    @{ Old = "Advanced Callouts"; New = "3.5 Advanced Callouts" },
    @{ Old = "Columns Layout"; New = "3.6 Columns Layout" },
    @{ Old = "Working with Package Functions"; New = "4 Working with Package Functions" },
    @{ Old = "Code Annotations"; New = "5 Code Annotations" },
    @{ Old = "Diagrams with Mermaid"; New = "6 Diagrams with Mermaid" },
    @{ Old = "Mathematical Notation"; New = "7 Mathematical Notation" },
    @{ Old = "Tips and Best Practices"; New = "8 Tips and Best Practices" },
    @{ Old = "Conclusion"; New = "9 Conclusion" },
    @{ Old = "Learn More"; New = "10 Learn More" },
    @{ Old = "References"; New = "11 References" }
)

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Range.Style.NameLocal
    if ($styleName -eq "Heading 2" -or $styleName -eq "Heading 3") {
        $text = $p.Range.Text.TrimEnd([char]13, [char]7)
        foreach ($r in $replacements) {
            if ($text -eq $r.Old) {
                $rng = $p.Range
                $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
                break
            }
        }
    }
}
